# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# for the ca9ac1fd-ed02-42d9-ad23-dfa79d5c9ffc file after a new handback
# round-trip was generated.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
# Row 3 corresponds to ca9ac1fd-ed02-42d9-ad23-dfa79d5c9ffc.md
# Column G = "Latest HO Xliff Generate Date"
$wsOverview.Range("G3").Value = "2016-08-15 09:02:59"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Row 3 corresponds to ca9ac1fd-ed02-42d9-ad23-dfa79d5c9ffc.md
# Column H = "Correspond Handoff Datetime"
$wsZhCn.Range("H3").Value = "2016-08-15 09:02:54"
# Column K = "Correspond Handback DateTime"
$wsZhCn.Range("K3").Value = "2016-08-15 09:03:25"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Row 3 corresponds to ca9ac1fd-ed02-42d9-ad23-dfa79d5c9ffc.md
# Column H = "Correspond Handoff Datetime"
$wsDeDe.Range("H3").Value = "2016-08-15 09:02:59"
# Column K = "Correspond Handback DateTime"
$wsDeDe.Range("K3").Value = "2016-08-15 09:03:32"
